$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (avg_glucose_level)
$ws.Range("B2").Value = 107.1061711711712
$ws.Range("C2").Value = 89.21600000000001
$ws.Range("D2").Value = 106.7146138807429
$ws.Range("E2").Value = 118.5287924528302
$ws.Range("F2").Value = 95.18405063291139
$ws.Range("G2").Value = 103.3499256275311

# Row 3 (age)
$ws.Range("B3").Value = 52.2972972972973
$ws.Range("C3").Value = 15.2
$ws.Range("D3").Value = 45.6236559139785
$ws.Range("E3").Value = 59.81132075471698
$ws.Range("F3").Value = 6.729789029535866
$ws.Range("G3").Value = 35.93241259910572

# Row 4 (bmi)
$ws.Range("B4").Value = 30.7972972972973
$ws.Range("C4").Value = 24.46
$ws.Range("D4").Value = 30.37429130009775
$ws.Range("E4").Value = 30.18339622641509
$ws.Range("F4").Value = 19.88354430379747
$ws.Range("G4").Value = 27.13970582552152
